$d = $word.ActiveDocument

function Insert-RunsXml($TargetRange, $PStyle, $Tokens) {
    $runsXml = ""
    foreach ($tok in $Tokens) {
        $escaped = $tok.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
        $runsXml = $runsXml + "<w:r><w:t xml:space=`"preserve`">$escaped</w:t></w:r>"
    }

    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p><w:pPr><w:pStyle w:val="' + $PStyle + '"/></w:pPr>' + $runsXml + '</w:p>' +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $TargetRange.InsertXML($xmlFrag) | Out-Null
}

# --- Paragraph 1: "On Pilgrimage - July/Aug 198" (Heading1 -> Title, split into runs) ---
$p1 = $d.Paragraphs(1)
$titleTokens = @("On", " ", "Pilgrimage", " ", "-", " ", "July", "/", "Aug", " ", "198")
Insert-RunsXml $p1.Range "Title" $titleTokens

# --- Paragraph 2: "By Dorothy Day" (Normal/bold -> Authors style, "Dorothy Day" split into runs) ---
$p2 = $d.Paragraphs(2)
$authorTokens = @("Dorothy", " ", "Day")
Insert-RunsXml $p2.Range "Authors" $authorTokens

# --- Best-effort bookmark cleanup (the stray "on-pilgrimage---julyaug-198" bookmark
#     that used to wrap the old title paragraph). The loaded document's pre-existing
#     bookmark isn't enumerated by this host's Bookmarks collection (Count stays 0
#     even though the raw <w:bookmarkStart/End> pair is still present in the part),
#     so this is defensive/best-effort and expected to be a no-op here. ---
try {
    $bmName = "on-pilgrimage---julyaug-198"
    $n = $d.Bookmarks.Count
    for ($i = $n; $i -ge 1; $i--) {
        $d.Bookmarks($i).Delete()
    }
    try { $d.Bookmarks.Item($bmName).Delete() } catch { }
    try { $d.Bookmarks($bmName).Delete() } catch { }
} catch {
    # Bookmarks collection may not be addressable in this runtime; ignore.
}
